$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell updates that can be written directly (Excel will not
# misinterpret these as numbers/dates: percentages with
# surrounding spaces, coin names, and URLs).
$plainUpdates = @{
    "D2" = '24.760.25'
    "E2" = '  +2.40%  '
    "D3" = '1.704.87'
    "E3" = '  +1.62%  '
    "E4" = '  -0.44%  '
    "E5" = '  +0.04%  '
    "E6" = '  -0.19%  '
    "E7" = '  +0.29%  '
    "E8" = '  +3.49%  '
    "E9" = '  -0.06%  '
    "E10" = '  +0.41%  '
    "E11" = '  +2.13%  '
    "E12" = '  -0.58%  '
    "E13" = '  +2.15%  '
    "E14" = '  +2.50%  '
    "E15" = '  +2.68%  '
    "D16" = '1.705.85'
    "E16" = '  +1.29%  '
    "E17" = '  +1.53%  '
    "E18" = '  +0.13%  '
    "E19" = '  -0.08%  '
    "E20" = '  +2.86%  '
    "E21" = '  +4.06%  '
    "E22" = '  +3.85%  '
    "E23" = '  +8.43%  '
    "D24" = '24.731.98'
    "E24" = '  +2.36%  '
    "E25" = '  +0.21%  '
    "E26" = '  +3.80%  '
    "E27" = '  +2.99%  '
    "E28" = '  -1.83%  '
    "E29" = '  +3.28%  '
    "D30" = '1.893.09'
    "E30" = '  +1.19%  '
    "E31" = '  +20.88%  '
    "E32" = '  +5.84%  '
    "E33" = '  +1.31%  '
    "B34" = 'WEMIXTOKEN'
    "C34" = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
    "E34" = '  +0.16%  '
    "B35" = 'Stellar'
    "C35" = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
    "E35" = '  +4.32%  '
    "B36" = 'Aptos'
    "C36" = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
    "E36" = '  +10.97%  '
    "E37" = '  +3.19%  '
    "E38" = '  +1.80%  '
    "E39" = '  +0.77%  '
    "B40" = 'VeChain'
    "C40" = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
    "E40" = '  +1.98%  '
    "B41" = 'Algorand'
    "C41" = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
    "E41" = '  +5.44%  '
    "E42" = '  +0.15%  '
    "E43" = '  +4.21%  '
    "E44" = '  -0.08%  '
    "E45" = '  +5.07%  '
    "E46" = '  +2.64%  '
    "E47" = '  +0.30%  '
    "E48" = '  +3.98%  '
    "E49" = '  +2.46%  '
    "E50" = '  +1.37%  '
    "E51" = '  +4.04%  '
}

foreach ($addr in $plainUpdates.Keys) {
    $ws.Range($addr).Value = $plainUpdates[$addr]
}

# Price values that look like plain numbers to Excel's type
# inference (e.g. "1.788", "309.00") must be forced to text so
# they are stored verbatim, matching the original inline-string
# cells. NumberFormat is reset back to the default afterwards so
# no stray cell style is left behind.
$textUpdates = @{
    "D4" = '1.001'
    "D5" = '309.00'
    "D6" = '0.9971'
    "D7" = '0.3739'
    "D8" = '49.17'
    "D9" = '0.3443'
    "D10" = '1.190'
    "D11" = '0.07461'
    "D12" = '0.9971'
    "D13" = '20.92'
    "D14" = '6.248'
    "D15" = '6.952'
    "D17" = '0.00001125'
    "D18" = '0.06717'
    "D19" = '0.9976'
    "D20" = '84.05'
    "D21" = '17.12'
    "D22" = '6.333'
    "D23" = '13.02'
    "D25" = '2.427'
    "D26" = '2.776'
    "D27" = '20.16'
    "D28" = '150.37'
    "D29" = '131.03'
    "D31" = '1.187'
    "D32" = '6.783'
    "D33" = '4.162'
    "D34" = '1.788'
    "D35" = '0.08819'
    "D36" = '13.69'
    "D37" = '5.533'
    "D38" = '0.06574'
    "D39" = '8.972'
    "D40" = '0.02385'
    "D41" = '0.2235'
    "D42" = '1.276'
    "D43" = '0.6429'
    "D44" = '0.9969'
    "D45" = '13.85'
    "D46" = '0.6108'
    "D48" = '2.113'
    "D49" = '129.99'
    "D50" = '0.07283'
    "D51" = '79.19'
}

foreach ($addr in $textUpdates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $textUpdates[$addr]
    $cell.Style = "Normal"
}

Write-Output "cryptos list updated"
